$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "print label" table, shrink the sheet back to just the
#     header + 3 sample rows that the new appsettings-driven export produces.
$ws.Rows("5:10").Delete()

$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# --- New header row (7 columns instead of the old 2) --------------------
# (leading "'" forces text + keeps the quote-prefix cell style the first two
#  header cells already had, same as the old "@@CTRLSUM@@"/"@@SHELF@@" ones)
$ws.Range("A1").Value = "'MHA    "
$ws.Range("B1").Value = "'Rack    "
$ws.Range("C1").Value = "X-Coor    "
$ws.Range("D1").Value = "Y-Coor    "
$ws.Range("E1").Value = "LocType    "
$ws.Range("F1").Value = "CheckSum    "
$ws.Range("G1").Value = "Zone    "

# --- Sample data rows, one per shelf level -------------------------------
$ws.Range("A2").Value = "BGN1 "
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "BP1"
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = "A "

$ws.Range("A3").Value = "BGN1 "
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "BP1"
$ws.Range("F3").Value = 92
$ws.Range("G3").Value = "A "

$ws.Range("A4").Value = "BGN1 "
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "BP1"
$ws.Range("F4").Value = 28
$ws.Range("G4").Value = "A "

# --- Move the selection / window to where the author left it ------------
$ws.Range("C7").Select()
$excel.ActiveWindow.Left = 2820
